# Updates symbol list (cryptos.xlsx) with refreshed prices / rankings.
# Cells D (Price) hold numeric-looking values but are stored as TEXT in this
# workbook, so we prefix them with a leading apostrophe to force Excel to
# keep them as text (preserving exact formatting / trailing zeros) instead
# of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BNB)
$ws.Range("D2").Value = '''247.38'

# Row 3 (OKB)
$ws.Range("D3").Value = '''22.00'

# Row 5 (Cronos)
$ws.Range("D5").Value = '''0.05775'

# Row 6 (GateToken)
$ws.Range("D6").Value = '''3.401'

# Row 7 (KuCoinToken)
$ws.Range("D7").Value = '''6.372'

# Row 8 (MXToken)
$ws.Range("D8").Value = '''0.8186'

# Row 9 (FTXToken)
$ws.Range("D9").Value = '''0.9694'
$ws.Range("E9").Value = '8FTXTokenFTT'

# Row 10 -> now "One" (rows 10-18 shift down from the old row above them)
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '''0.01120'
$ws.Range("E10").Value = '9OneONEBestin24h'

# Row 11 -> now "WazirX"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1430'
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12 -> now "MandalaExchangeToken"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07452'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

# Row 13 -> now "LiechtensteinCryptoassetsExchange"
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '''0.03136'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

# Row 14 -> now "BitrueCoin"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.02998'
$ws.Range("E14").Value = '13BitrueCoinBTR'

# Row 15 -> now "MCDex"
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '''4.148'
$ws.Range("E15").Value = '14MCDexMCB'

# Row 16 -> now "BitMartToken"
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = '''0.09416'
$ws.Range("E16").Value = '15BitMartTokenBMX'

# Row 17 -> now "BitForexToken"
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '''0.001607'
$ws.Range("E17").Value = '16BitForexTokenBF'

# Row 18 -> now "CoinExToken"
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '''0.04822'
$ws.Range("E18").Value = '17CoinExTokenCET'

# Remaining scattered price refreshes
$ws.Range("D20").Value = '''0.004135'
$ws.Range("D21").Value = '''0.0009948'
$ws.Range("D22").Value = '''0.0001499'
$ws.Range("D23").Value = '''3.767'
$ws.Range("D25").Value = '''0.3258'
$ws.Range("D26").Value = '''0.1261'
$ws.Range("D27").Value = '''0.0003998'
$ws.Range("D40").Value = '''0.03896'
$ws.Range("D41").Value = '''0.006458'
$ws.Range("D43").Value = '''0.002999'
$ws.Range("D44").Value = '''0.006244'
$ws.Range("D45").Value = '''0.00005594'
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("D47").Value = '''0.3799'
$ws.Range("D49").Value = '''0.00002099'
$ws.Range("D50").Value = '''0.01010'
